# Applies the "cryptos" price/volume refresh described by the commit diff.
# Rows 22/23, 32/33 and 45/46 were re-ranked (their Coin/Link/Price/Volume swapped),
# all other rows simply got refreshed Price (column D) and Volume(1h) (column E) values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unicode subscript-3 character used inside one PEPE price (0.0<sub3>0788).
$sub3 = [char]0x2083

$ws.Range("D2").Value = "64.233.92"
$ws.Range("E2").Value = "  -0.51%  "

$ws.Range("D3").Value = "3.159.55"
$ws.Range("E3").Value = "  -0.54%  "

$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").Value = "'611.81"
$ws.Range("E5").Value = "  +2.02%  "

$ws.Range("D6").Value = "'147.84"
$ws.Range("E6").Value = "  -2.78%  "

$ws.Range("E7").Value = "  +0.11%  "

$ws.Range("D8").Value = "3.151.65"
$ws.Range("E8").Value = "  -0.84%  "

$ws.Range("D9").Value = "'0.526"
$ws.Range("E9").Value = "  -0.44%  "

$ws.Range("E10").Value = "  -0.91%  "

$ws.Range("E11").Value = "  -1.67%  "

$ws.Range("D12").Value = "'0.472"
$ws.Range("E12").Value = "  -0.64%  "

$ws.Range("D13").Value = "'0.0000259"
$ws.Range("E13").Value = "  +0.41%  "

$ws.Range("D14").Value = "'35.59"
$ws.Range("E14").Value = "  -3.62%  "

$ws.Range("D15").Value = "3.673.97"
$ws.Range("E15").Value = "  -0.86%  "

$ws.Range("E16").Value = "  +3.01%  "

$ws.Range("D17").Value = "64.210.89"
$ws.Range("E17").Value = "  -0.67%  "

$ws.Range("D18").Value = "3.156.35"
$ws.Range("E18").Value = "  -0.89%  "

$ws.Range("D19").Value = "'6.91"
$ws.Range("E19").Value = "  -1.69%  "

$ws.Range("D20").Value = "'479.29"
$ws.Range("E20").Value = "  -0.21%  "

$ws.Range("D21").Value = "'14.74"
$ws.Range("E21").Value = "  -0.47%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'8.09"
$ws.Range("E22").Value = "  +4.44%  "

$ws.Range("B23").Value = "Polygon"
$ws.Range("C23").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D23").Value = "'0.717"
$ws.Range("E23").Value = "  +0.20%  "

$ws.Range("D24").Value = "'13.76"
$ws.Range("E24").Value = "  -0.89%  "

$ws.Range("D25").Value = "'83.81"
$ws.Range("E25").Value = "  -1.12%  "

$ws.Range("E26").Value = "  +0.11%  "

$ws.Range("E27").Value = "  -3.07%  "

$ws.Range("D28").Value = "'8.54"
$ws.Range("E28").Value = "  -1.01%  "

$ws.Range("D29").Value = "'7.18"
$ws.Range("E29").Value = "  +3.14%  "

$ws.Range("D30").Value = "'0.121"
$ws.Range("E30").Value = "  -0.81%  "

$ws.Range("E31").Value = "  -6.10%  "

$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  -0.02%  "

$ws.Range("B33").Value = "Stacks"
$ws.Range("C33").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D33").Value = "'2.72"
$ws.Range("E33").Value = "  -0.34%  "

$ws.Range("D34").Value = "'26.39"
$ws.Range("E34").Value = "  -1.73%  "

$ws.Range("E35").Value = "  +1.72%  "

$ws.Range("D36").Value = "0.0${sub3}0788"
$ws.Range("E36").Value = "  +6.91%  "

$ws.Range("D37").Value = "'6.02"
$ws.Range("E37").Value = "  -1.61%  "

$ws.Range("D38").Value = "'52.97"
$ws.Range("E38").Value = "  -3.05%  "

$ws.Range("D39").Value = "'3.17"
$ws.Range("E39").Value = "  -2.28%  "

$ws.Range("D40").Value = "'463.36"
$ws.Range("E40").Value = "  +1.33%  "

$ws.Range("D41").Value = "'0.0400"
$ws.Range("E41").Value = "  -0.52%  "

$ws.Range("E42").Value = "  -4.40%  "

$ws.Range("E43").Value = "  -1.48%  "

$ws.Range("D44").Value = "2.870.55"
$ws.Range("E44").Value = "  -0.19%  "

$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").Value = "'0.270"
$ws.Range("E45").Value = "  -1.90%  "

$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").Value = "'2.31"
$ws.Range("E46").Value = "  -4.59%  "

$ws.Range("E47").Value = "  +4.34%  "

$ws.Range("E48").Value = "  -2.64%  "

$ws.Range("E49").Value = "  -0.02%  "

$ws.Range("E50").Value = "  -1.33%  "

$ws.Range("D51").Value = "'119.23"
$ws.Range("E51").Value = "  -0.80%  "
